$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells (order matters: it controls the order new shared
# strings are appended to the workbook's shared-string table, mirroring
# the author's edit order: AE1, Z1, AA1, AB1, X1, W1) ---
$ws.Range("AE1").Value = "z_interest"
$ws.Range("Z1").Value  = "x_size"
$ws.Range("AA1").Value = "y_size"
$ws.Range("AB1").Value = "z_size"
$ws.Range("X1").Value  = "z1_ind"
$ws.Range("W1").Value  = "z0_ind"

# --- Remove the now-unwanted derived (dx_ind / dy_ind / dz_ind) formula
# columns S, V, Y for every data row (2-14) ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("S$r").ClearContents()
    $ws.Range("V$r").ClearContents()
    $ws.Range("Y$r").ClearContents()
}

# --- Update the active selection / view to match the saved state ---
$ws.Range("S1").Select() | Out-Null
